$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 846, shifting the existing rows 846-947
# down to 848-949 (dimension grows from A1:R947 to A1:R949).
$ws.Range("A846:R847").EntireRow.Insert()

# Populate the two newly inserted rows (846 and 847) with the new
# weekly records. All the "constant" columns for this sheet
# (A, B, C, E, F, G, H, N, Q, R) keep the same values as every other
# row in the block.

# Row 846: "Primera"
$ws.Cells.Item(846, 1).Value = 9
$ws.Cells.Item(846, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(846, 3).Value = "Metropolitana"
$ws.Cells.Item(846, 4).Value = 45142
$ws.Cells.Item(846, 5).Value = 13
$ws.Cells.Item(846, 6).Value = 100114014
$ws.Cells.Item(846, 7).Value = "Betarraga"
$ws.Cells.Item(846, 8).Value = "Sin especificar"
$ws.Cells.Item(846, 9).Value = "Primera"
$ws.Cells.Item(846, 10).Value = 7900
$ws.Cells.Item(846, 11).Value = 90
$ws.Cells.Item(846, 12).Value = 100
$ws.Cells.Item(846, 13).Value = 95
$ws.Cells.Item(846, 14).Value = "$/unidad"
$ws.Cells.Item(846, 15).Value = "Región Metropolitana"
$ws.Cells.Item(846, 16).Value = 95
$ws.Cells.Item(846, 17).Value = 1
$ws.Cells.Item(846, 18).Value = "Hortaliza"

# Row 847: "Segunda"
$ws.Cells.Item(847, 1).Value = 9
$ws.Cells.Item(847, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(847, 3).Value = "Metropolitana"
$ws.Cells.Item(847, 4).Value = 45142
$ws.Cells.Item(847, 5).Value = 13
$ws.Cells.Item(847, 6).Value = 100114014
$ws.Cells.Item(847, 7).Value = "Betarraga"
$ws.Cells.Item(847, 8).Value = "Sin especificar"
$ws.Cells.Item(847, 9).Value = "Segunda"
$ws.Cells.Item(847, 10).Value = 5200
$ws.Cells.Item(847, 11).Value = 70
$ws.Cells.Item(847, 12).Value = 70
$ws.Cells.Item(847, 13).Value = 70
$ws.Cells.Item(847, 14).Value = "$/unidad"
$ws.Cells.Item(847, 15).Value = "Región Metropolitana"
$ws.Cells.Item(847, 16).Value = 70
$ws.Cells.Item(847, 17).Value = 1
$ws.Cells.Item(847, 18).Value = "Hortaliza"
